# Add a new "Long term loans" data column (G) to both worksheets.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1 (2)")
$ws2 = $wb.Worksheets.Item("Sheet1")

# --- Header ---------------------------------------------------------
$ws1.Range("G1").Value = "Long term loans"
$ws2.Range("G1").Value = "Long term loans"

# --- "Sheet1 (2)" (reverse chronological, rows 2-56) -----------------
$data1 = @(17.82, 18.89, 22.99, 18.82, 14.56, 28.08, 24.37, 33.63, 29.55, 30.32, 29.62, 26.41, 22.78, 4.94, -14.82, -18.33, -21.77, -15.34, -45.95, -27.57, -21.3, -15.13, 6.75, 3.47, 17.98, 20.7, 22.69, 23.24, 7.51, -24.72, -30.44, -2.88, -17.71, -14.32, -13.04, -11.47, -1.29, -3.41, 6.54, 17.7, 20.75, 15.78, 18.28, 36.93, 32.26, 31.96, 29.98, 29.1, 27.34, 16.4, 6.82, 16.32, 18.81, 23.92, 18.32)

$r = 2
foreach ($v in $data1) {
    $ws1.Cells.Item($r, 7).Value = $v
    $r++
}

# --- "Sheet1" (chronological, rows 2-73) -----------------------------
$data2 = @(8.37, -0.92, 19.77, -12.88, -0.52, 2.75, 4.54, 6.64, 13.98, 20.49, 23.46, 18.32, 23.92, 18.81, 16.32, 6.82, 16.4, 27.34, 29.1, 29.98, 31.96, 32.26, 36.93, 18.28, 15.78, 20.75, 17.7, 6.54, -3.41, -1.29, -11.47, -13.04, -14.32, -17.71, -2.88, -30.44, -24.72, 7.51, 23.24, 22.69, 20.7, 17.98, 3.47, 6.75, -15.13, -21.3, -27.57, -45.95, -15.34, -21.77, -18.33, -14.82, 4.94, 22.78, 26.41, 29.62, 30.32, 29.55, 33.63, 24.37, 28.08, 14.56, 18.82, 22.99, 18.89, 17.82, 21.94, 31.54, 6.72, 0.72, -17.08, -1.06)

$r = 2
foreach ($v in $data2) {
    $ws2.Cells.Item($r, 7).Value = $v
    $r++
}

# --- Selections (match the authored workbook's last active cells) ---
$ws2.Activate()
$ws2.Range("H17").Select()

$ws1.Activate()
$ws1.Range("I15").Select()
